$wb = $excel.ActiveWorkbook

$wsAddCustomer = $wb.Worksheets.Item("AddCustomerTest")

# "AddCustomerTest" becomes the active/selected tab (was "test_suite" before)
$wsAddCustomer.Activate()

# The active cell/selection on "AddCustomerTest" moves from E3 to E2
$wsAddCustomer.Range("E2").Select()

# The "runmode" column (E2:E4) flips from "Y" to "N"
$wsAddCustomer.Range("E2").Value = "N"
$wsAddCustomer.Range("E3").Value = "N"
$wsAddCustomer.Range("E4").Value = "N"

$wb.Save()
